$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates reflecting refreshed crypto price/volume data from the source feed.
# NumberFormat is forced to text ("@") before assignment so that values such as
# "1.000", "2.320" or "0.00001125" are preserved exactly as strings instead of
# being normalized by Excel into numeric values (e.g. 1, 2.32, 1.125E-05).
$updates = @(
    @{ Cell = 'D2'; Value = '28.252.90' }
    @{ Cell = 'E2'; Value = '  +0.03%  ' }
    @{ Cell = 'D3'; Value = '1.908.95' }
    @{ Cell = 'E3'; Value = '  +2.05%  ' }
    @{ Cell = 'D4'; Value = '1.000' }
    @{ Cell = 'E4'; Value = '  -0.08%  ' }
    @{ Cell = 'D5'; Value = '314.51' }
    @{ Cell = 'E5'; Value = '  +0.96%  ' }
    @{ Cell = 'E6'; Value = '  -0.03%  ' }
    @{ Cell = 'D7'; Value = '0.5067' }
    @{ Cell = 'E7'; Value = '  -0.03%  ' }
    @{ Cell = 'D8'; Value = '0.3928' }
    @{ Cell = 'E8'; Value = '  +0.25%  ' }
    @{ Cell = 'E9'; Value = '  -3.27%  ' }
    @{ Cell = 'D10'; Value = '1.142' }
    @{ Cell = 'E10'; Value = '  +0.03%  ' }
    @{ Cell = 'D11'; Value = '41.83' }
    @{ Cell = 'E11'; Value = '  +2.23%  ' }
    @{ Cell = 'D12'; Value = '6.403' }
    @{ Cell = 'E12'; Value = '  -1.70%  ' }
    @{ Cell = 'B13'; Value = 'Solana' }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol' }
    @{ Cell = 'D13'; Value = '20.89' }
    @{ Cell = 'E13'; Value = '  -0.56%  ' }
    @{ Cell = 'B14'; Value = 'WrappedEther' }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' }
    @{ Cell = 'D14'; Value = '1.914.63' }
    @{ Cell = 'E14'; Value = '  +2.16%  ' }
    @{ Cell = 'D15'; Value = '7.311' }
    @{ Cell = 'E15'; Value = '  -1.80%  ' }
    @{ Cell = 'D16'; Value = '1.001' }
    @{ Cell = 'E16'; Value = '  -0.05%  ' }
    @{ Cell = 'D17'; Value = '0.00001125' }
    @{ Cell = 'E17'; Value = '  -0.40%  ' }
    @{ Cell = 'D18'; Value = '92.78' }
    @{ Cell = 'E18'; Value = '  -0.25%  ' }
    @{ Cell = 'D19'; Value = '0.06603' }
    @{ Cell = 'E19'; Value = '  +0.30%  ' }
    @{ Cell = 'D20'; Value = '17.98' }
    @{ Cell = 'E20'; Value = '  +2.28%  ' }
    @{ Cell = 'D21'; Value = '0.9991' }
    @{ Cell = 'E21'; Value = '  -0.13%  ' }
    @{ Cell = 'D22'; Value = '6.199' }
    @{ Cell = 'E22'; Value = '  +0.44%  ' }
    @{ Cell = 'D23'; Value = '28.300.09' }
    @{ Cell = 'E23'; Value = '  +0.02%  ' }
    @{ Cell = 'D24'; Value = '11.43' }
    @{ Cell = 'E24'; Value = '  +0.65%  ' }
    @{ Cell = 'D25'; Value = '2.320' }
    @{ Cell = 'E25'; Value = '  +1.31%  ' }
    @{ Cell = 'B26'; Value = 'LidoDAOToken' }
    @{ Cell = 'C26'; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo' }
    @{ Cell = 'D26'; Value = '2.592' }
    @{ Cell = 'E26'; Value = '  +1.45%  ' }
    @{ Cell = 'B27'; Value = 'WrappedliquidstakedEther2.0' }
    @{ Cell = 'C27'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth' }
    @{ Cell = 'D27'; Value = '2.130.83' }
    @{ Cell = 'E27'; Value = '  +1.96%  ' }
    @{ Cell = 'D28'; Value = '21.06' }
    @{ Cell = 'E28'; Value = '  -0.80%  ' }
    @{ Cell = 'D29'; Value = '157.79' }
    @{ Cell = 'E29'; Value = '  -0.40%  ' }
    @{ Cell = 'D30'; Value = '127.32' }
    @{ Cell = 'E30'; Value = '  -0.31%  ' }
    @{ Cell = 'D31'; Value = '1.102' }
    @{ Cell = 'E31'; Value = '  +2.92%  ' }
    @{ Cell = 'D32'; Value = '0.1073' }
    @{ Cell = 'E32'; Value = '  +0.95%  ' }
    @{ Cell = 'D33'; Value = '5.641' }
    @{ Cell = 'E33'; Value = '  +0.02%  ' }
    @{ Cell = 'D34'; Value = '3.617' }
    @{ Cell = 'E34'; Value = '  -0.22%  ' }
    @{ Cell = 'D35'; Value = '9.693' }
    @{ Cell = 'E35'; Value = '  +1.27%  ' }
    @{ Cell = 'D36'; Value = '0.06668' }
    @{ Cell = 'E36'; Value = '  -0.87%  ' }
    @{ Cell = 'D37'; Value = '0.02416' }
    @{ Cell = 'E37'; Value = '  +1.25%  ' }
    @{ Cell = 'D38'; Value = '1.250' }
    @{ Cell = 'E38'; Value = '  +0.77%  ' }
    @{ Cell = 'D39'; Value = '0.2190' }
    @{ Cell = 'E39'; Value = '  -0.02%  ' }
    @{ Cell = 'D40'; Value = '1.276' }
    @{ Cell = 'E40'; Value = '  +7.52%  ' }
    @{ Cell = 'D41'; Value = '0.6422' }
    @{ Cell = 'E41'; Value = '  +0.69%  ' }
    @{ Cell = 'D42'; Value = '5.013' }
    @{ Cell = 'E42'; Value = '  +0.61%  ' }
    @{ Cell = 'D43'; Value = '11.49' }
    @{ Cell = 'E43'; Value = '  -0.11%  ' }
    @{ Cell = 'D44'; Value = '1.000' }
    @{ Cell = 'E44'; Value = '  -0.02%  ' }
    @{ Cell = 'D45'; Value = '13.33' }
    @{ Cell = 'E45'; Value = '  -1.88%  ' }
    @{ Cell = 'D46'; Value = '0.6008' }
    @{ Cell = 'E46'; Value = '  -0.19%  ' }
    @{ Cell = 'D47'; Value = '3.720' }
    @{ Cell = 'E47'; Value = '  +1.58%  ' }
    @{ Cell = 'E48'; Value = '  +1.19%  ' }
    @{ Cell = 'D49'; Value = '2.020' }
    @{ Cell = 'E49'; Value = '  +0.98%  ' }
    @{ Cell = 'D50'; Value = '122.93' }
    @{ Cell = 'E50'; Value = '  -0.95%  ' }
    @{ Cell = 'E51'; Value = '  -0.93%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $range.NumberFormat = '@'
    $range.Value = $u.Value
}
